# Convert date of birth to datetime isoformat
#
# The "birthdate" column (G) previously stored dates as plain integers
# typed like MMDDYYYY (e.g. 10201980 for 1980-10-20). This replaces those
# values with proper Excel date serial numbers and applies a date number
# format so they round-trip as real dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mary Watson - birthdate 10/20/1980
$ws.Range("G2").Value = 29514
# John Newman - birthdate 10/20/1992
$ws.Range("G3").Value = 33897

# Apply a date number format (maps to the builtin "mm-dd-yy" / numFmtId 14)
$ws.Range("G2").NumberFormat = "mm-dd-yy"

# Re-use the same style for G3 (copy/paste-format) instead of re-deriving a
# number format, so both cells share one style entry
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)

# Widen column G so the formatted dates are fully visible
$ws.Columns.Item(7).ColumnWidth = 21.74

# Reflect the final selection left after making the edit
$ws.Range("K3").Select()
